# "add evidence to B1~B6"
# Fill in the placeholder TxHash rows on sheets B5 and B6 with the real
# evidence hash values, and add a new sheet "B7" with its own TxHash
# evidence pair, following the same layout as the existing B1/B2/B5/B6
# sheets (row1 = "TxHash" header, row2/row3 = hash values).

$wb = $excel.ActiveWorkbook

# --- B1: used below as the formatting source for B7's header cell ---
$wsB1 = $wb.Worksheets.Item("B1")

# --- B2: just a selection/cursor move, no content change ---
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Activate()
$wsB2.Range("A3").Select()

# --- B5: replace placeholder text with the real hash evidence ---
$wsB5 = $wb.Worksheets.Item("B5")
$wsB5.Range("A2").Value = "19680C762BF4B581EC9BC45E2AD3CEFC6713E0BE7D6C75A73A2360DAE144F41E"
$wsB5.Range("A3").Value = "D4D24BBE5D0FC4554F8DDC455D2D5C2E19D801A7B3DBE77D1185B6C6085C662F"
$wsB5.Activate()
$wsB5.Range("A3").Select()

# --- B6: replace placeholder text with the real hash evidence ---
$wsB6 = $wb.Worksheets.Item("B6")
$wsB6.Range("A2").Value = "2015195D628E78D1709DB59A530ED4A27262A617996891516A95C2DB787F8C39"
$wsB6.Range("A3").Value = "24E99CA3ABC00B6C77D4B1CA314EED1BE5140F7256B51B5ABBBDF8A64504255A"
$wsB6.Activate()
$wsB6.Range("A1").Select()

# --- B7: brand-new sheet appended after B6, same TxHash layout ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsB7 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsB7.Name = "B7"

# Match the "TxHash" header look used on the other evidence sheets
# (bold header row with a light-grey fill) by copying A1's formatting
# from an existing evidence sheet.
$wsB1.Range("A1").Copy()
$wsB7.Range("A1").PasteSpecial(-4122)

$wsB7.Range("A1").Value = "TxHash"
$wsB7.Range("A2").Value = "F7137C060DABAC746C84A289067904FE76C630B3752DFFEA9660423D1B9CB8DD"
$wsB7.Range("A3").Value = "845EF268B95B1585F009036154C549562F60828333D476F34B206AC08E1347FE"
$wsB7.Activate()
$wsB7.Range("A3").Select()
